# Apply "change int to float" edit:
# - Add a new note in cell C1 of each of the three sheets explaining that the
#   column must be filled in when a skill is actually used.
# - Update the active selection on each sheet to reflect where the user left
#   off editing.

$wb = $excel.ActiveWorkbook

$note = "（若技能需要使用，该列一定要进行填写）"

$ws1 = $wb.Worksheets.Item("Side1")
$ws1.Range("C1").Value = $note
$ws1.Range("C9").Select()

$ws2 = $wb.Worksheets.Item("Side2")
$ws2.Range("C1").Value = $note
$ws2.Range("C1").Select()

$ws3 = $wb.Worksheets.Item("Side3")
$ws3.Range("C1").Value = $note
$ws3.Range("C8").Select()

$ws1.Activate()
